$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# New translated-text values that become new shared strings (si 97..117,
# since uniqueCount grows from 97 to 118).
# -----------------------------------------------------------------------
$s97 = ' Hey, hey, hey! I\''m going after\noutlaws from the Outlaw Notice Board today!'
$s98 = ' Let\''s look good training! Hey,\nhey, hey!'
$s99 = 'SCRIPT/G01P03A/um1202.ssb'
$s100 = ' Эй, эй, эй! Сегодня я буду\nохотиться за негодяями с Доски\nРазыскиваемых!'
$s101 = ' Пора хорошенько потренироваться!\nЭй, эй, эй!'
$s102 = ' Üê, üê, üê! Òåãïäîÿ ÿ áôäô\nïöïóéóûòÿ èà îåãïäÿÿíé ò Äïòëé\nÑàèúòëéâàåíúö!'
$s103 = ' Ðïñà öïñïšåîûëï ðïóñåîéñïâàóûòÿ!\nÜê, üê, üê!'
$s104 = 'SCRIPT/T01P01A/um1304.ssb'
$s105 = ' Let\''s catch that thieving [CS:N]Grovyle[CR],\nhey, hey!'
$s106 = ' Мы поймаем этого похитителя\n[CS:N]Гровайла[CR], эй, эй!'
$s107 = ' Íú ðïêíàåí üóïãï ðïöéóéóåìÿ\n[CS:N]Ãñïâàêìà[CR], üê, üê!'
$s108 = ' Hey, hey! Not knowing what\''s\nbeen happening…'
$s109 = ' Hey, hey! That\''s stressful, I say.'
$s110 = ' Well, all we can do is keep doing\nour usual work, hey, hey.'
$s111 = 'SCRIPT/G01P03A/um1610.ssb'
$s112 = ' Эй, эй! Незнание происходящего...'
$s113 = ' Эй, эй! Должен сказать, это меня\nнапрягает.'
$s114 = ' Ну, всё что мы можем сделать,\nэто работать как и прежде, эй, эй.'
$s115 = ' Üê, üê! Îåèîàîéå ðñïéòöïäÿþåãï...'
$s116 = ' Üê, üê! Äïìçåî òëàèàóû, üóï íåîÿ\nîàðñÿãàåó.'
$s117 = ' Îô, âòæ œóï íú íïçåí òäåìàóû,\nüóï ñàáïóàóû ëàë é ðñåçäå, üê, üê.'

# =========================================================================
# Row 28 (existing last row of the table): only its cell style changes, from
# the "normal" style (s=4 numeric/text col A-B, s=5 text col C-E) to the
# "closing border" style (s=9 / s=10) because a new block (rows 29-34) is
# appended after it. Values and row height (57.6) are unchanged.
# =========================================================================
$origA28 = $ws.Range("A28").Value2
$origB28 = $ws.Range("B28").Value2
$origC28 = $ws.Range("C28").Value2
$origD28 = $ws.Range("D28").Value2
$origE28 = $ws.Range("E28").Value2
$ws.Range("A23:E23").Copy($ws.Range("A28:E28"))
$ws.Range("A28").Value = $origA28
$ws.Range("B28").Value = $origB28
$ws.Range("C28").Value = $origC28
$ws.Range("D28").Value = $origD28
$ws.Range("E28").Value = $origE28

# =========================================================================
# New rows 29-34, using the same visual style patterns already present in
# the sheet (copy style from the matching existing row, then overwrite the
# cell values/row height for the new content).
# =========================================================================

# --- Row 29: style like row 26 (s=4 / s=5), height 43.2 ---
$ws.Range("A26:E26").Copy($ws.Range("A29:E29"))
$ws.Range("A29").Value = $s99
$ws.Range("B29").Value = 482
$ws.Range("C29").Value = $s97
$ws.Range("D29").Value = $s100
$ws.Range("E29").Value = $s102
$ws.Rows.Item(29).RowHeight = 43.2

# --- Row 30: style like row 20 (s=7 / s=8, empty A cell), height 21.6 ---
$ws.Range("A20:E20").Copy($ws.Range("A30:E30"))
$ws.Range("B30").Value = 485
$ws.Range("C30").Value = $s98
$ws.Range("D30").Value = $s101
$ws.Range("E30").Value = $s103
$ws.Rows.Item(30).RowHeight = 21.6

# --- Row 31: style like row 23 (s=9 / s=10), height 43.2 ---
$ws.Range("A23:E23").Copy($ws.Range("A31:E31"))
$ws.Range("A31").Value = $s104
$ws.Range("B31").Value = 463
$ws.Range("C31").Value = $s105
$ws.Range("D31").Value = $s106
$ws.Range("E31").Value = $s107
$ws.Rows.Item(31).RowHeight = 43.2

# --- Row 32: style like row 26 (s=4 / s=5), height 43.2 ---
$ws.Range("A26:E26").Copy($ws.Range("A32:E32"))
$ws.Range("A32").Value = $s111
$ws.Range("B32").Value = 438
$ws.Range("C32").Value = $s108
$ws.Range("D32").Value = $s112
$ws.Range("E32").Value = $s115
$ws.Rows.Item(32).RowHeight = 43.2

# --- Row 33: style like row 26 cols B-E (s=4 / s=5), NO cell in col A, height 21.6 ---
$ws.Range("B26:E26").Copy($ws.Range("B33:E33"))
$ws.Range("B33").Value = 441
$ws.Range("C33").Value = $s109
$ws.Range("D33").Value = $s113
$ws.Range("E33").Value = $s116
$ws.Rows.Item(33).RowHeight = 21.6

# --- Row 34: style like row 26 cols B-E (s=4 / s=5), NO cell in col A, height 21.6 ---
$ws.Range("B26:E26").Copy($ws.Range("B34:E34"))
$ws.Range("B34").Value = 444
$ws.Range("C34").Value = $s110
$ws.Range("D34").Value = $s114
$ws.Range("E34").Value = $s117
$ws.Rows.Item(34).RowHeight = 21.6

# =========================================================================
# View-state: selection/active cell moves to E34 (best-effort; the engine
# does not persist scroll position (topLeftCell), only the selection).
# =========================================================================
$null = $ws.Range("E34").Select()

Write-Output "Applied Korfish sheet update: row 28 restyled, rows 29-34 added."